$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F273").Value = 34582
$ws.Range("G273").Value = 1822
$ws.Range("F282").Value = 47127
$ws.Range("G282").Value = 2840
$ws.Range("F283").Value = 16932
$ws.Range("G283").Value = 1003
$ws.Range("F291").Value = 14954
$ws.Range("G291").Value = 486
$ws.Range("F292").Value = 81774
$ws.Range("G292").Value = 7253
$ws.Range("F294").Value = 92105
$ws.Range("G294").Value = 4864
$ws.Range("F317").Value = 61674
$ws.Range("F318").Value = 49865
$ws.Range("F319").Value = 41308
$ws.Range("F320").Value = 68485
$ws.Range("G320").Value = 3122
$ws.Range("F321").Value = 92083
$ws.Range("G321").Value = 2815
$ws.Range("F324").Value = 234297
$ws.Range("G324").Value = 2686
$ws.Range("F325").Value = 759622
$ws.Range("G325").Value = 6452
$ws.Range("F326").Value = 431076
$ws.Range("G326").Value = 3815
$ws.Range("F327").Value = 236911
$ws.Range("G327").Value = 2873
$ws.Range("F328").Value = 180216
$ws.Range("G328").Value = 2638
$ws.Range("F329").Value = 89206
$ws.Range("F332").Value = 432645
$ws.Range("G332").Value = 4238
$ws.Range("F333").Value = 262702
$ws.Range("G333").Value = 2835
$ws.Range("F334").Value = 203646
$ws.Range("G334").Value = 3394
$ws.Range("F335").Value = 128731
$ws.Range("G335").Value = 2882
$ws.Range("F336").Value = 100689
$ws.Range("G336").Value = 3199
$ws.Range("F337").Value = 102394
$ws.Range("G337").Value = 2903
$ws.Range("F338").Value = 217307
$ws.Range("G338").Value = 3058
$ws.Range("F339").Value = 638530
$ws.Range("G339").Value = 5430
$ws.Range("F340").Value = 378458
$ws.Range("G340").Value = 3236
$ws.Range("F341").Value = 296095
$ws.Range("G341").Value = 3665
$ws.Range("F342").Value = 171268
$ws.Range("G342").Value = 2891
$ws.Range("F343").Value = 126117
$ws.Range("G343").Value = 2812
$ws.Range("F344").Value = 129749
$ws.Range("G344").Value = 2387
$ws.Range("F345").Value = 274447
$ws.Range("G345").Value = 3128
$ws.Range("F346").Value = 626572
$ws.Range("G346").Value = 4427
$ws.Range("F347").Value = 325419
$ws.Range("G347").Value = 2770
$ws.Range("F348").Value = 220383
$ws.Range("G348").Value = 2996
$ws.Range("F349").Value = 156507
$ws.Range("G349").Value = 2643

# Add new row 350
$ws.Range("A350").Value = 44244
$ws.Range("A350").NumberFormat = "yyyy-mm-dd"
$ws.Range("B350").Value = 285419
$ws.Range("C350").Value = 10848
$ws.Range("D350").Value = 2555
$ws.Range("E350").Value = 6271
$ws.Range("F350").Value = 105334
$ws.Range("G350").Value = 2248
